# "replaced all reading times with average"
#
# For each of the three per-tree-count sheets (10_trees, 500_trees, 1600_trees):
#   - add a new row 88 with D88 = "average"
#   - E88 on 10_trees computes the grand AVERAGE of the "data loading time"
#     column (E2:E86) across all three sheets
#   - E88 on 500_trees / 1600_trees simply mirrors '10_trees'!E88
#   - every K column formula (SUM(E{r},G{r},H{r})) is rewritten to use the
#     new average cell instead of the row's own E value:
#     SUM($E$88,G{r},H{r})
# The L column (PRODUCT(K,0.001)) and everything else is left untouched and
# recalculates naturally because it depends on K.

$wb = $excel.ActiveWorkbook

$sheetNames = @("10_trees", "500_trees", "1600_trees")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("D88").Value = "average"

    if ($name -eq "10_trees") {
        $ws.Range("E88").Formula = "=AVERAGE(E2:E86,'500_trees'!E2:E86,'1600_trees'!E2:E86)"
    } else {
        $ws.Range("E88").Formula = "='10_trees'!E88"
    }

    $ws.Range("K2").Formula = "=SUM(`$E`$88,G2,H2)"
    $ws.Range("K3:K86").Formula = "=SUM(`$E`$88,G3,H3)"
}
